$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.786833763122559
$ws.Range("B1").Value = 4.698285579681396
$ws.Range("C1").Value = 2.137770891189575
$ws.Range("D1").Value = 1.521677136421204
$ws.Range("E1").Value = 1.304303646087646
